# Automatische test-sync: 2025-07-29 21:51:50
#
# Appends the newest "Retour / Terugbetaling" mail log entry (Testmail #11)
# to the Logs sheet, bumps its matching Dashboard category tally, and
# extends the conditional formatting + chart series ranges so the new row
# is covered.

$wb = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

# --- Logs!A13:J13 -----------------------------------------------------
$logs.Cells.Item(13, 1).Value = "Mijn retour is nog steeds niet verwerkt."
$logs.Cells.Item(13, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(13, 3).Value = "Testmail #11: Mijn retour is nog steeds niet verwerkt."
$logs.Cells.Item(13, 4).Value = "Retour / Terugbetaling"
$logs.Cells.Item(13, 5).Value = "Beste klant,`nBedankt voor je bericht. Het spijt me te horen dat je retour nog niet is verwerkt. Om dit verder te kunnen onderzoeken, heb ik wat meer informatie van je nodig. Zou je alsjeblieft je ordernummer en/of het trackingnummer van de retourzending kunnen doorgeven? Hiermee kunnen we precies nakijken wat er gaande is en je zo snel mogelijk helpen.`nMet vriendelijke groet,`n[Naam]`nKlantenservice Team"
$logs.Cells.Item(13, 6).Value = "2025-07-29 21:51:07"
$logs.Cells.Item(13, 7).Value = "Ja"
$logs.Cells.Item(13, 8).Value = "Nee"
$logs.Cells.Item(13, 9).Value = "Ja"
$logs.Cells.Item(13, 10).Value = "Nee"

# the multi-line Antwoord text triggers an auto row-height bump; put the
# row back to its (implicit) default height like every other data row
$logs.Rows.Item(13).AutoFit()

# --- grow the conditional formatting ranges to include row 13 --------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range($col + "2:" + $col + "12")
    $newRange = $logs.Range($col + "2:" + $col + "13")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Dashboard!A6:B6 (new category tally row) -------------------------
$dashboard.Cells.Item(6, 1).Value = "Retour / Terugbetaling"
$dashboard.Cells.Item(6, 2).Value = 1

# --- extend the bar chart's category/value series to row 6 -----------
$chartObj = $dashboard.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.XValues = '''Dashboard''!$A$2:$A$6'
$series.Values = '''Dashboard''!$B$2:$B$6'

Write-Output "Logs/Dashboard synced through row 13/6"
